$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update "Valor Mora" total summary cell (E11)
$ws.Range("E11").Value = 227760

# Shift the "Periodo Mora" (E16:E19) values forward one period:
# old: 2507, 2506, 2505, 2504 -> new: 2505, 2506, 2507, 2508
$ws.Range("E16").Value = "2505"
$ws.Range("E17").Value = "2506"
$ws.Range("E18").Value = "2507"
$ws.Range("E19").Value = "2508"

# Update "Valor Mora" for the last row (F19) to match the new period's value
$ws.Range("F19").Value = 56940
